$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.444740653038025
$ws.Range("B1").Value = 3.31209659576416
$ws.Range("C1").Value = 4.187763214111328
$ws.Range("D1").Value = 1.943168520927429
$ws.Range("E1").Value = 1.156373739242554
